$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 757 ("今日辛いこと、明日は平気だろう" post), which shifts
# all subsequent rows (758-797) up by one.
$ws.Rows.Item(757).Delete()
